# Update cryptos list prices/volume percentages (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells stay text (some values are plain-looking numbers
# like "226.63" that Excel would otherwise auto-convert to floating point).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.019.09"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.015.77"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.63"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.607"
$ws.Range("E6").Value = "  -1.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.78"
$ws.Range("E8").Value = "  -3.90%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.316.86"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.24"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.26"
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.742"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.011.28"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.985.32"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.05"
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.82"
$ws.Range("E20").Value = "  -1.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.98"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("E25").Value = "  -5.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.08"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.125"
$ws.Range("E28").Value = "  -3.16%  "
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "18.72"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("E35").Value = "  -5.72%  "
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -4.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.38"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0216"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.476.24"
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "95.08"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.54"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("E45").Value = "  -4.98%  "
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.23"
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.92"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.204.88"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.22"
$ws.Range("E51").Value = "  -2.53%  "
